$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date serial numbers that were all bumped
# from 45186 (2023-09-17) to 45188 (2023-09-19) for rows 2 through 121.
$ws.Range("C2:C121").Value = 45188
